$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rangeA = $ws.Range("B$r1`:AC$r1")
    $rangeB = $ws.Range("B$r2`:AC$r2")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

# Swap the data (columns B..AC) between the following row pairs.
# Column A (the running id) is left untouched.
Swap-Rows 129 130
Swap-Rows 131 132
Swap-Rows 145 146
Swap-Rows 149 151

# Append a new fixture row (row 152) with no result yet (future match).
$ws.Range("A151").Copy()
$ws.Range("A152").PasteSpecial(-4122)
$ws.Range("E151").Copy()
$ws.Range("E152").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 7867520
$ws.Range("C152").Value = "South Korea K3 League"
$ws.Range("D152").Value = "South Korea K3 League"
$ws.Range("E152").Value = 45387.125
$ws.Range("F152").Value = "Daejeon Korail"
$ws.Range("G152").Value = "Gimhae City"
$ws.Range("K152").Value = 2.8
$ws.Range("L152").Value = 3
$ws.Range("M152").Value = 2.375
$ws.Range("N152").Value = 3.1
$ws.Range("O152").Value = 3
$ws.Range("P152").Value = 2.3
$ws.Range("Q152").Value = 0.25
$ws.Range("R152").Value = 1.8
$ws.Range("S152").Value = 2
$ws.Range("T152").Value = 2.25
$ws.Range("U152").Value = 2
$ws.Range("V152").Value = 1.8
$ws.Range("W152").Value = 0
$ws.Range("X152").Value = 0
$ws.Range("Y152").Value = 0
$ws.Range("Z152").Value = 0
$ws.Range("AA152").Value = 0
